$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per latest scrape.
# D-column values are written via NumberFormat="@" + Style reset so
# numeric-looking strings (e.g. "242.77") stay text, matching the source data.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '97.366.35'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.43%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.140.31'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.62%  '

$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.77'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.38%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '612.77'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.10%  '

$ws.Range("E7").Value = '  +2.28%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.385'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.51%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.02%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.137.56'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.59%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.783'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.43%  '

$ws.Range("E12").Value = '  +0.28%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '97.059.89'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.39%  '

$ws.Range("E14").Value = '  -1.70%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.50'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.13%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '34.17'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.72%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.727.18'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.48%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.138.44'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.29%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.57'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.33%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '521.12'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +18.06%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.66'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.39%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.69'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.31%  '

$ws.Range("E23").Value = '  -4.33%  '

$ws.Range("E24").Value = '  -3.15%  '

$ws.Range("E25").Value = '  -2.36%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '88.85'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.44%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.68'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -6.71%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.308.10'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.45%  '

$ws.Range("E29").Value = '  +0.03%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.241'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.25%  '

$ws.Range("E31").Value = '  -1.96%  '

$ws.Range("E32").Value = '  +0.95%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.996'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.27%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '9.05'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.30%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '26.72'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.58%  '

$ws.Range("E36").Value = '  -4.13%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.43'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -8.91%  '

$ws.Range("E38").Value = '  -0.71%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '24.24'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.05%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '478.90'
$ws.Range("D40").Style = "Normal"

$ws.Range("E41").Value = '  +1.52%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.23'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.05%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.60'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -9.90%  '

$ws.Range("E44").Value = '  +0.02%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.17'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.53%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '160.69'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.06%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.51'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.82%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '44.18'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.55%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.997'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.11%  '

# Row 47/48: ARBITRUM and Stacks swapped positions (rank reorder) with updated values
$ws.Range("B47").Value = 'Stacks'
$ws.Range("C47").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.94'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.11%  '

$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.705'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.38%  '
